# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# 1. Ativacao date text update (B8, C8)
Set-TextValue $ws.Range("B8") '01/01/2022'
Set-TextValue $ws.Range("C8") '01/01/2022'

# 2. Docentes responsaveis - update existing row 13
Set-TextValue $ws.Range("B13") '5840897 - Clodoaldo Saron'
Set-TextValue $ws.Range("C13") '5840897 - Clodoaldo Saron'

# 3. Insert a new row at 14 for the second docente
$ws.Rows(14).Insert()
Set-TextValue $ws.Range("B14") '1033242 - Fábio Herbst Florenzano'
Set-TextValue $ws.Range("C14") '1033242 - Fábio Herbst Florenzano'

# 4. Update Metodo (now row 20), Criterio (row 21), Norma de recuperacao (row 22), Bibliografia (row 23)
Set-TextValue $ws.Range("B20") 'A avaliação será feita por meio de Provas Escritas, Estudos de Casos e Desenvolvimento de Projetos, sendo necessário utilizar pelo menos dois critérios de avaliação diferentes.'
Set-TextValue $ws.Range("C20") 'A avaliação será feita por meio de Provas Escritas, Estudos de Casos e Desenvolvimento de Projetos, sendo necessário utilizar pelo menos dois critérios de avaliação diferentes.'

Set-TextValue $ws.Range("B21") 'A Nota final (NF) será calculada da seguinte maneira: NF = (P+EC+Projetos)/3'
Set-TextValue $ws.Range("C21") 'A Nota final (NF) será calculada da seguinte maneira: NF = (P+EC+Projetos)/3'

Set-TextValue $ws.Range("B22") 'Não consta recuperação'
Set-TextValue $ws.Range("C22") 'Não consta recuperação'

Set-TextValue $ws.Range("B23") '1. J. Margolis. Engineering Plastics Handbook. McGraw-Hill Professional, 2005. 2. Nigel Mills. Plastics - Microstructure and Engineering Applications. Butterworth-Heineman, 2005. 3. Walter Michaeli, TEcnologia dos Plasticos. Ed. Blucher 4. Hélio Wiebeck, Júlio Harada. Plásticos de Engenharia - Tecnologia e Aplicações. São Paulo: Editora Artliber, 2005. 5. E. B. Mano, L. C. Mendes. Identificação de Plásticos, Borrachas e Fibras. São Paulo: Editora Edgard Blucher, 2000. 6. Marcelo Rabello. Aditivação de Polímeros. São Paulo: Editora Artliber, 2004. 7. Jan C.J. Bart. Additives in Polymers. New York: John Wiley & Sons, 2005. 8. Marino Xanthos. Functional Fillers for Plastics. Wiley-VCH Verlag GmbH, 2005. 9. Silvio Manrich. Processamento de Termoplásticos. Editora Artliber, 2005. 10. G.H. Michler, F.J. Baltá-Calleja. Mechanical Properties of Polymers Based on Nanostructure and Morphology. Boca Raton: CRC Press, 2005. 11. A. M. Piva, H. Wiebeck. Reciclagem do P. São Paulo: Editora Artliber". Manas Chanda, ,Salil K. Roy  Plastics Fabrication and Recycling'
Set-TextValue $ws.Range("C23") '1. J. Margolis. Engineering Plastics Handbook. McGraw-Hill Professional, 2005. 2. Nigel Mills. Plastics - Microstructure and Engineering Applications. Butterworth-Heineman, 2005. 3. Walter Michaeli, TEcnologia dos Plasticos. Ed. Blucher 4. Hélio Wiebeck, Júlio Harada. Plásticos de Engenharia - Tecnologia e Aplicações. São Paulo: Editora Artliber, 2005. 5. E. B. Mano, L. C. Mendes. Identificação de Plásticos, Borrachas e Fibras. São Paulo: Editora Edgard Blucher, 2000. 6. Marcelo Rabello. Aditivação de Polímeros. São Paulo: Editora Artliber, 2004. 7. Jan C.J. Bart. Additives in Polymers. New York: John Wiley & Sons, 2005. 8. Marino Xanthos. Functional Fillers for Plastics. Wiley-VCH Verlag GmbH, 2005. 9. Silvio Manrich. Processamento de Termoplásticos. Editora Artliber, 2005. 10. G.H. Michler, F.J. Baltá-Calleja. Mechanical Properties of Polymers Based on Nanostructure and Morphology. Boca Raton: CRC Press, 2005. 11. A. M. Piva, H. Wiebeck. Reciclagem do P. São Paulo: Editora Artliber". Manas Chanda, ,Salil K. Roy  Plastics Fabrication and Recycling'

$excel.CutCopyMode = $false
